$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Switch the workbook's default (Normal) font to the CJK font seen in the
# target file (宋体 / SimSun) - this is what the author's Excel did when it
# re-saved the sheet with its regional/default-font settings.
$normal = $wb.Styles.Item("Normal")
$normal.Font.Name = "宋体"

# Add the new value in column D, continuing the 1,2,3 sequence with a 4.
$ws.Range("D1").Value = 4

# Excel leaves the active cell/selection on the cell that was just edited.
$ws.Range("D1").Select()

# Match the printer/page setup recorded for the sheet.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
